$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.854.57"
$ws.Range("E2").Value = "  -0.62%  "

$ws.Range("D3").Value = "1.636.45"
$ws.Range("E3").Value = "  -0.18%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -1.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.22"
$ws.Range("E5").Value = "  -0.81%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5030"
$ws.Range("E6").Value = "  +0.65%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.004"
$ws.Range("E7").Value = "  -1.11%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2559"
$ws.Range("E8").Value = "  -0.83%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06370"
$ws.Range("E9").Value = "  -0.86%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.39"
$ws.Range("E10").Value = "  -0.66%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07782"
$ws.Range("E11").Value = "  +0.17%  "

$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.244"
$ws.Range("E12").Value = "  -0.42%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.635.33"
$ws.Range("E13").Value = "  -0.34%  "

$ws.Range("D14").Value = "1.861.94"
$ws.Range("E14").Value = "  -0.10%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5402"
$ws.Range("E15").Value = "  -1.07%  "

$ws.Range("D16").Value = "0.0₅7876"
$ws.Range("E16").Value = "  -0.88%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.35"
$ws.Range("E17").Value = "  +1.00%  "

$ws.Range("D18").Value = "25.883.40"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.005"
$ws.Range("E19").Value = "  -0.96%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "195.49"
$ws.Range("E20").Value = "  -4.81%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.362"
$ws.Range("E21").Value = "  +1.14%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.884"
$ws.Range("E22").Value = "  -1.32%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.951"
$ws.Range("E23").Value = "  -0.44%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.005"
$ws.Range("E24").Value = "  -0.89%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.891"
$ws.Range("E25").Value = "  -4.00%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "139.84"
$ws.Range("E26").Value = "  -1.47%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1132"
$ws.Range("E27").Value = "  -1.90%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.796"
$ws.Range("E28").Value = "  -0.11%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.63"
$ws.Range("E29").Value = "  -0.92%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.237"
$ws.Range("E30").Value = "  -0.49%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04834"
$ws.Range("E31").Value = "  -3.91%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.238"
$ws.Range("E32").Value = "  -0.99%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.164"
$ws.Range("E33").Value = "  -1.36%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.527"
$ws.Range("E34").Value = "  -1.09%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.363"
$ws.Range("E35").Value = "  +0.48%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.601"
$ws.Range("E36").Value = "  -0.35%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.8841"
$ws.Range("E37").Value = "  -0.90%  "

$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5500"
$ws.Range("E38").Value = "  -2.76%  "

$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "1.124.02"
$ws.Range("E39").Value = "  -0.63%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01557"
$ws.Range("E40").Value = "  -0.49%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.004"
$ws.Range("E41").Value = "  -1.09%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.656"
$ws.Range("E42").Value = "  +0.24%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8114"
$ws.Range("E43").Value = "  -0.85%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.23"
$ws.Range("E44").Value = "  -0.67%  "

$ws.Range("D45").Value = "0.0₈122"
$ws.Range("E45").Value = "  +10.11%  "

$ws.Range("D46").Value = "1.773.43"
$ws.Range("E46").Value = "  +0.00%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4523"
$ws.Range("E47").Value = "  -0.60%  "

$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.005"
$ws.Range("E48").Value = "  -0.79%  "

$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "55.01"
$ws.Range("E49").Value = "  +0.30%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05037"
$ws.Range("E50").Value = "  -0.18%  "

$ws.Range("E51").Value = "  -0.49%  "
